$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 17.99273490905762
$ws.Range("D2").Value = 26377

$ws.Range("C3").Value = 16.96491241455078
$ws.Range("D3").Value = 330

$ws.Range("C4").Value = 16.50500297546387
$ws.Range("D4").Value = 242

$ws.Range("C5").Value = 16.51978492736816
$ws.Range("D5").Value = 242

$ws.Range("C6").Value = 16.5712833404541
$ws.Range("D6").Value = 242
